$d = $word.ActiveDocument

$d.Content.Find.Execute(
    " The consolidation reported sufficient staff at that time to meet this portion of the requirements.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)
